# Update the NATMI LR-pairs sheet (Icosl-Ctla4) with refreshed TPM-based
# values: rows 5:7 are removed entirely, and rows 2:4 are rewritten with the
# recomputed stats (row 2 now pairs ECs -> MuSCs instead of ECs -> ECs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three trailing rows (old rows 5,6,7) - this also shrinks the
# used range / dimension down to A1:T4 and shifts nothing else up since
# they are the last rows.
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2: ECs -> Icosl/Ctla4 -> MuSCs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icosl"
$ws.Range("C2").Value = "Ctla4"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6721510000000001
$ws.Range("H2").Value = 2.016453
$ws.Range("I2").Value = 0.03553601427691108
$ws.Range("J2").Value = 0.03553601427691108
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01388066666666667
$ws.Range("N2").Value = 0.041642
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.009329903980666666
$ws.Range("R2").Value = 0.083969135826
$ws.Range("S2").Value = 0.03553601427691108
$ws.Range("T2").Value = 0.03553601427691108

# Row 3: FAPs -> Icosl/Ctla4 -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Icosl"
$ws.Range("C3").Value = "Ctla4"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 15.25501433333334
$ws.Range("H3").Value = 45.76504300000001
$ws.Range("I3").Value = 0.8065187839396453
$ws.Range("J3").Value = 0.8065187839396453
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01388066666666667
$ws.Range("N3").Value = 0.041642
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.2117497689562222
$ws.Range("R3").Value = 1.905747920606
$ws.Range("S3").Value = 0.8065187839396453
$ws.Range("T3").Value = 0.8065187839396453

# Row 4: MuSCs -> Icosl/Ctla4 -> MuSCs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Icosl"
$ws.Range("C4").Value = "Ctla4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.987477
$ws.Range("H4").Value = 8.962431
$ws.Range("I4").Value = 0.1579452017834437
$ws.Range("J4").Value = 0.1579452017834437
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01388066666666667
$ws.Range("N4").Value = 0.041642
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.04146817241133333
$ws.Range("R4").Value = 0.373213551702
$ws.Range("S4").Value = 0.1579452017834437
$ws.Range("T4").Value = 0.1579452017834437
